$wb = $excel.ActiveWorkbook

# Add "Computer" sheet after the last existing sheet ("Device")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsComputer = $wb.Worksheets.Add($null, $lastSheet)
$wsComputer.Name = "Computer"
$wsComputer.Range("A1").Value = "os"
$wsComputer.Range("B1").Value = "osversion"
$wsComputer.Range("C1").Value = "location"
$wsComputer.Range("D1").Value = "vm"

# Add "Harddisk" sheet after "Computer"
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHarddisk = $wb.Worksheets.Add($null, $lastSheet2)
$wsHarddisk.Name = "Harddisk"
$wsHarddisk.Range("A1").Value = "size"
$wsHarddisk.Range("B1").Value = "formfactor"
$wsHarddisk.Range("C1").Value = "rpm"
$wsHarddisk.Range("D1").Value = "computer"
